$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook gained one new weekly data record ("Fruta / hortaliza, semanal").
# In the canonical OOXML, this shows up as a brand-new row 99, with every
# existing row from the old 99 through 210 shifted down by one (to 100-211).
# Insert a new row at position 99 to push the existing data down, then
# populate it with the new record's values.

$ws.Rows.Item(99).Insert()

$ws.Cells.Item(99, 1).Value2  = 5
$ws.Cells.Item(99, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(99, 3).Value2  = "Maule"
$ws.Cells.Item(99, 4).Value2  = 44539
$ws.Cells.Item(99, 5).Value2  = 7
$ws.Cells.Item(99, 6).Value2  = 100112003
$ws.Cells.Item(99, 7).Value2  = "Ajo"
$ws.Cells.Item(99, 8).Value2  = "Chino"
$ws.Cells.Item(99, 9).Value2  = "Primera"
$ws.Cells.Item(99, 10).Value2 = 250
$ws.Cells.Item(99, 11).Value2 = 18000
$ws.Cells.Item(99, 12).Value2 = 18000
$ws.Cells.Item(99, 13).Value2 = 18000
$ws.Cells.Item(99, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(99, 15).Value2 = "China"
$ws.Cells.Item(99, 16).Value2 = 1800
$ws.Cells.Item(99, 17).Value2 = 10
$ws.Cells.Item(99, 18).Value2 = "Hortaliza"
